# "Implementation de la bd plus debug"
# Adds two days of work (rows 15 & 16) to the "Iteration #3" sheet's
# work-log table, bumps the first day's hours, and moves the on-sheet
# selection to reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #3")
$ws.Activate()

# Day 1 (2018-04-23) already existed — its hours grew from 3 to 4.
$ws.Range("C14").Value = 4

# Day 2 (2018-04-24): date, task description, hours worked.
$ws.Range("A15").Value = 43214
$ws.Range("A15").NumberFormat = "d-mmm"
$ws.Range("B15").Value = "Suite implémentation de la BD dans l'app"
$ws.Range("C15").Value = 5

# Day 3 (2018-04-25): same task continued, different hours.
$ws.Range("A16").Value = 43215
$ws.Range("A16").NumberFormat = "d-mmm"
$ws.Range("B16").Value = "Suite implémentation de la BD dans l'app"
$ws.Range("C16").Value = 2

# Leave the user's selection on C17, just past the rows just filled in.
$ws.Range("C17").Select() | Out-Null
